$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new time-sheet entry (date, hours worked, blank time-formatted cell)
$ws.Range("A4").Value = 42918
$ws.Range("A4").NumberFormat = "MM/DD/YY"
$ws.Range("B4").Value = 0.5
$ws.Range("C4").NumberFormat = "HH:MM:SS\ AM/PM"

# Row 5: new time-sheet entry (date, blank hours cell, time value)
$ws.Range("A5").Value = 42918
$ws.Range("A5").NumberFormat = "MM/DD/YY"
$ws.Range("C5").Value = 0.388888888888889
$ws.Range("C5").NumberFormat = "HH:MM:SS\ AM/PM"

$ws.Range("C5").Select()
